# Add 2022-Q3 data.
#
# The existing "2022-Q2" sheet's numbers actually belong to the new
# "2022-Q3" period (values got revised), and its old/original numbers move
# to a brand-new "2022-Q2" sheet inserted right after it. "2022-Q1" and
# "2021-Q4" stay as-is (just shift position). The "总计" (summary) sheet
# gets a new row for 2022-Q3 inserted after its header.

$wb = $excel.ActiveWorkbook

# --- Worksheets ---------------------------------------------------------

# Duplicate the current "2022-Q2" sheet (placed immediately after it) so we
# keep its original numbers under the "2022-Q2" label, while the original
# sheet object gets relabeled "2022-Q3" and updated with the new figures.
$origQ2 = $wb.Worksheets.Item("2022-Q2")
$origQ2.Copy($null, $origQ2)
$newQ2 = $wb.Worksheets.Item("2022-Q2 (2)")

# Free up the "2022-Q2" name before claiming it on the copy.
$origQ2.Name = "2022-Q3"
$newQ2.Name = "2022-Q2"

# Update the "2022-Q3" sheet with the revised figures (kept as text, same
# as the rest of the column, not converted to numbers).
$q3 = $wb.Worksheets.Item("2022-Q3")
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q3.Range($addr).NumberFormat = "@"
}
$q3.Range("D2").Value = "0.45"
$q3.Range("E2").Value = "90.87"
$q3.Range("F2").Value = "1.90"
$q3.Range("G2").Value = "0.0086"
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q3.Range($addr).Style = "Normal"
}

# --- Summary ("总计") sheet ----------------------------------------------

$summary = $wb.Worksheets.Item(1)

# Shift rows 2-4 down to rows 3-5, then write the new 2022-Q3 row into row 2.
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.01

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 8
$summary.Range("D4").Value = 1.66

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.01

$summary.Range("B2").Value = "2022-Q3"

# New row 5 needs column A's formatting (bordered/centered header style)
# carried down, same as A2:A4.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
